$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: 2020-08-24
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "'2020-08-24"
$ws.Cells.Item(12, 2).Value = 3009
$ws.Cells.Item(12, 3).Value = 3621
$ws.Cells.Item(12, 4).Value = 612
$ws.Cells.Item(12, 5).Value = 4692.5
$ws.Cells.Item(12, 6).Value = 3759.5
$ws.Cells.Item(12, 7).Value = -933.0000000000001
$ws.Cells.Item(12, 8).Value = 5927
$ws.Cells.Item(12, 9).Value = 5047.5
$ws.Cells.Item(12, 10).Value = -879.5000000000002
$ws.Cells.Item(12, 11).Value = 504
$ws.Cells.Item(12, 12).Value = 265
$ws.Cells.Item(12, 13).Value = -239
$ws.Cells.Item(12, 14).Value = 14077.4
$ws.Cells.Item(12, 15).Value = 12526.4
$ws.Cells.Item(12, 16).Value = -1551.000000000001
$ws.Cells.Item(12, 17).Value = 7962.100000000001
$ws.Cells.Item(12, 18).Value = 7545.2
$ws.Cells.Item(12, 19).Value = -416.9000000000015
$ws.Cells.Item(12, 20).Value = 18953
$ws.Cells.Item(12, 21).Value = 19739
$ws.Cells.Item(12, 22).Value = 786
$ws.Cells.Item(12, 23).Value = 44296.9
$ws.Cells.Item(12, 24).Value = 46634.89999999999
$ws.Cells.Item(12, 25).Value = 2338.000000000002
$ws.Cells.Item(12, 26).Value = 53316.4
$ws.Cells.Item(12, 27).Value = 56225.39999999999
$ws.Cells.Item(12, 28).Value = 2909.000000000003
$ws.Cells.Item(12, 29).Value = 3113
$ws.Cells.Item(12, 30).Value = 1677
$ws.Cells.Item(12, 31).Value = -1436
$ws.Cells.Item(12, 32).Value = 1151.1
$ws.Cells.Item(12, 33).Value = 2978.2
$ws.Cells.Item(12, 34).Value = 748.9000000000001
$ws.Cells.Item(12, 35).Value = 1273.1
$ws.Cells.Item(12, 36).Value = 3427.2
$ws.Cells.Item(12, 37).Value = 901.9000000000001

# Row 13: 2020-08-25
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "'2020-08-25"
$ws.Cells.Item(13, 2).Value = 6607
$ws.Cells.Item(13, 3).Value = 7644
$ws.Cells.Item(13, 4).Value = 1037
$ws.Cells.Item(13, 5).Value = 6277.5
$ws.Cells.Item(13, 6).Value = 6762.1
$ws.Cells.Item(13, 7).Value = 484.5999999999997
$ws.Cells.Item(13, 8).Value = 9292
$ws.Cells.Item(13, 9).Value = 9757.6
$ws.Cells.Item(13, 10).Value = 465.5999999999997
$ws.Cells.Item(13, 11).Value = 1321
$ws.Cells.Item(13, 12).Value = 1930
$ws.Cells.Item(13, 13).Value = 705
$ws.Cells.Item(13, 14).Value = 19962.5
$ws.Cells.Item(13, 15).Value = 19401.8
$ws.Cells.Item(13, 16).Value = -560.6999999999998
$ws.Cells.Item(13, 17).Value = 14212.2
$ws.Cells.Item(13, 18).Value = 14587.1
$ws.Cells.Item(13, 19).Value = 374.8999999999993
$ws.Cells.Item(13, 20).Value = 31117
$ws.Cells.Item(13, 21).Value = 33654
$ws.Cells.Item(13, 22).Value = 2537
$ws.Cells.Item(13, 23).Value = 67795
$ws.Cells.Item(13, 24).Value = 68818.5
$ws.Cells.Item(13, 25).Value = 1023.500000000003
$ws.Cells.Item(13, 26).Value = 82518.49999999999
$ws.Cells.Item(13, 27).Value = 85215
$ws.Cells.Item(13, 28).Value = 2696.5
$ws.Cells.Item(13, 29).Value = 7249
$ws.Cells.Item(13, 30).Value = 1979
$ws.Cells.Item(13, 31).Value = -5157
$ws.Cells.Item(13, 32).Value = 3578.5
$ws.Cells.Item(13, 33).Value = 5134.2
$ws.Cells.Item(13, 34).Value = 1349.7
$ws.Cells.Item(13, 35).Value = 6214.5
$ws.Cells.Item(13, 36).Value = 5927.7
$ws.Cells.Item(13, 37).Value = -603.3000000000004

# Row 14: 2020-08-26
$ws.Range("A11").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "'2020-08-26"
$ws.Cells.Item(14, 2).Value = 3815
$ws.Cells.Item(14, 3).Value = 4421
$ws.Cells.Item(14, 4).Value = 343
$ws.Cells.Item(14, 5).Value = 4386.6
$ws.Cells.Item(14, 6).Value = 4012.2
$ws.Cells.Item(14, 7).Value = -374.3999999999997
$ws.Cells.Item(14, 8).Value = 5754.1
$ws.Cells.Item(14, 9).Value = 5714.2
$ws.Cells.Item(14, 10).Value = -39.90000000000009
$ws.Cells.Item(14, 11).Value = 1124
$ws.Cells.Item(14, 12).Value = 862
$ws.Cells.Item(14, 13).Value = -262
$ws.Cells.Item(14, 14).Value = 12739
$ws.Cells.Item(14, 15).Value = 13869.8
$ws.Cells.Item(14, 16).Value = 1130.800000000001
$ws.Cells.Item(14, 17).Value = 7692.799999999999
$ws.Cells.Item(14, 18).Value = 8144.9
$ws.Cells.Item(14, 19).Value = 452.1000000000008
$ws.Cells.Item(14, 20).Value = 21172
$ws.Cells.Item(14, 21).Value = 21122
$ws.Cells.Item(14, 22).Value = 732
$ws.Cells.Item(14, 23).Value = 46042.5
$ws.Cells.Item(14, 24).Value = 48598.10000000001
$ws.Cells.Item(14, 25).Value = 3344.100000000002
$ws.Cells.Item(14, 26).Value = 56115.5
$ws.Cells.Item(14, 27).Value = 58906.60000000001
$ws.Cells.Item(14, 28).Value = 3633.600000000002
$ws.Cells.Item(14, 29).Value = 2061
$ws.Cells.Item(14, 30).Value = 1912
$ws.Cells.Item(14, 31).Value = -731
$ws.Cells.Item(14, 32).Value = 4675.4
$ws.Cells.Item(14, 33).Value = 1924
$ws.Cells.Item(14, 34).Value = -2606.5
$ws.Cells.Item(14, 35).Value = 5509.9
$ws.Cells.Item(14, 36).Value = 2482.9
$ws.Cells.Item(14, 37).Value = -3009.5

# Row 15: 2020-08-27
$ws.Range("A11").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "'2020-08-27"
$ws.Cells.Item(15, 2).Value = 5631
$ws.Cells.Item(15, 3).Value = 4154
$ws.Cells.Item(15, 4).Value = -1477
$ws.Cells.Item(15, 5).Value = 3417.6
$ws.Cells.Item(15, 6).Value = 2800.4
$ws.Cells.Item(15, 7).Value = -617.2
$ws.Cells.Item(15, 8).Value = 5523.6
$ws.Cells.Item(15, 9).Value = 3981.9
$ws.Cells.Item(15, 10).Value = -1541.7
$ws.Cells.Item(15, 11).Value = 747
$ws.Cells.Item(15, 12).Value = 822
$ws.Cells.Item(15, 13).Value = 75
$ws.Cells.Item(15, 14).Value = 12618.6
$ws.Cells.Item(15, 15).Value = 10700.4
$ws.Cells.Item(15, 16).Value = -1918.2
$ws.Cells.Item(15, 17).Value = 7544.3
$ws.Cells.Item(15, 18).Value = 4787
$ws.Cells.Item(15, 19).Value = -2757.3
$ws.Cells.Item(15, 20).Value = 19506
$ws.Cells.Item(15, 21).Value = 18949
$ws.Cells.Item(15, 22).Value = -924
$ws.Cells.Item(15, 23).Value = 44783.1
$ws.Cells.Item(15, 24).Value = 41218.1
$ws.Cells.Item(15, 25).Value = -3564.999999999998
$ws.Cells.Item(15, 26).Value = 51872.1
$ws.Cells.Item(15, 27).Value = 50059.10000000001
$ws.Cells.Item(15, 28).Value = -2468.999999999998
$ws.Cells.Item(15, 29).Value = 2344
$ws.Cells.Item(15, 30).Value = 1472
$ws.Cells.Item(15, 31).Value = -872
$ws.Cells.Item(15, 32).Value = 723.3
$ws.Cells.Item(15, 33).Value = 1135.7
$ws.Cells.Item(15, 34).Value = -39.79999999999995
$ws.Cells.Item(15, 35).Value = 1051.8
$ws.Cells.Item(15, 36).Value = 1555.2
$ws.Cells.Item(15, 37).Value = -57.79999999999995

# Row 16: 2020-08-28
$ws.Range("A11").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "'2020-08-28"
$ws.Cells.Item(16, 2).Value = 13386
$ws.Cells.Item(16, 3).Value = 10829
$ws.Cells.Item(16, 4).Value = -2557
$ws.Cells.Item(16, 5).Value = 11313
$ws.Cells.Item(16, 6).Value = 12856.6
$ws.Cells.Item(16, 7).Value = 1543.600000000001
$ws.Cells.Item(16, 8).Value = 16082.5
$ws.Cells.Item(16, 9).Value = 16462.1
$ws.Cells.Item(16, 10).Value = 379.6000000000004
$ws.Cells.Item(16, 11).Value = 1369
$ws.Cells.Item(16, 12).Value = 1873
$ws.Cells.Item(16, 13).Value = 504
$ws.Cells.Item(16, 14).Value = 33058.2
$ws.Cells.Item(16, 15).Value = 35604.2
$ws.Cells.Item(16, 16).Value = 1731.8
$ws.Cells.Item(16, 17).Value = 19371.4
$ws.Cells.Item(16, 18).Value = 20040.4
$ws.Cells.Item(16, 19).Value = 668.9999999999995
$ws.Cells.Item(16, 20).Value = 61699
$ws.Cells.Item(16, 21).Value = 66591
$ws.Cells.Item(16, 22).Value = 4892
$ws.Cells.Item(16, 23).Value = 128039.7
$ws.Cells.Item(16, 24).Value = 131441.7
$ws.Cells.Item(16, 25).Value = 3402.00000000001
$ws.Cells.Item(16, 26).Value = 157939.7
$ws.Cells.Item(16, 27).Value = 163283.7
$ws.Cells.Item(16, 28).Value = 5344.000000000011
$ws.Cells.Item(16, 29).Value = 9230
$ws.Cells.Item(16, 30).Value = 6630
$ws.Cells.Item(16, 31).Value = -2600
$ws.Cells.Item(16, 32).Value = 7226
$ws.Cells.Item(16, 33).Value = 4572.7
$ws.Cells.Item(16, 34).Value = -1044.1
$ws.Cells.Item(16, 35).Value = 10714.5
$ws.Cells.Item(16, 36).Value = 5946.7
$ws.Cells.Item(16, 37).Value = -869.1000000000004
